$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rough")
$ws.Range("E5:E16").HorizontalAlignment = -4131
